{"js": "const body = context.document.body;\nconst results = body.search(\"Fall 2025\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'Fall 2025' in the document body.\");\n}\n\nresults.items[0].insertText(\"Spring 2026\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Fall 2025\"\n$find.Replacement.Text = \"Spring 2026\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n) | Out-Null\n"}
